$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Delete the last 4 data rows (rows 6-9) so only rows 1-5 remain.
$ws.Rows("6:9").Delete()

# Row 2: ECs / Cirbp / Trem1 / Resolving-Mac
$ws.Range("A2").Value = "ECs"
$ws.Range("D2").Value = "Resolving-Mac"
$ws.Range("G2").Value = 18.496333
$ws.Range("H2").Value = 55.488999
$ws.Range("I2").Value = 0.3902867652967038
$ws.Range("J2").Value = 0.3902867652967038
$ws.Range("M2").Value = 0.09428199999999999
$ws.Range("N2").Value = 0.282846
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 1.743871267906
$ws.Range("R2").Value = 15.694841411154
$ws.Range("S2").Value = 0.3902867652967038
$ws.Range("T2").Value = 0.3902867652967038

# Row 3: FAPs / Cirbp / Trem1 / Resolving-Mac
$ws.Range("A3").Value = "FAPs"
$ws.Range("D3").Value = "Resolving-Mac"
$ws.Range("G3").Value = 12.34551533333333
$ws.Range("H3").Value = 37.036546
$ws.Range("I3").Value = 0.260499810712076
$ws.Range("J3").Value = 0.260499810712076
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.09428199999999999
$ws.Range("N3").Value = 0.282846
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 1.163959876657333
$ws.Range("R3").Value = 10.475638889916
$ws.Range("S3").Value = 0.260499810712076
$ws.Range("T3").Value = 0.260499810712076

# Row 4: MuSCs / Cirbp / Trem1 / Resolving-Mac
$ws.Range("A4").Value = "MuSCs"
$ws.Range("D4").Value = "Resolving-Mac"
$ws.Range("G4").Value = 11.24015833333333
$ws.Range("H4").Value = 33.720475
$ws.Range("I4").Value = 0.2371759330532953
$ws.Range("J4").Value = 0.2371759330532953
$ws.Range("M4").Value = 0.09428199999999999
$ws.Range("N4").Value = 0.282846
$ws.Range("O4").Value = 1
$ws.Range("P4").Value = 1
$ws.Range("Q4").Value = 1.059744607983333
$ws.Range("R4").Value = 9.537701471849999
$ws.Range("S4").Value = 0.2371759330532953
$ws.Range("T4").Value = 0.2371759330532953

# Row 5: Resolving-Mac / Cirbp / Trem1 / Resolving-Mac
$ws.Range("A5").Value = "Resolving-Mac"
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("G5").Value = 5.309641333333333
$ws.Range("H5").Value = 15.928924
$ws.Range("I5").Value = 0.112037490937925
$ws.Range("J5").Value = 0.1120374909379251
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.09428199999999999
$ws.Range("N5").Value = 0.282846
$ws.Range("O5").Value = 1
$ws.Range("P5").Value = 1
$ws.Range("Q5").Value = 0.5006036041893333
$ws.Range("R5").Value = 4.505432437704
$ws.Range("S5").Value = 0.112037490937925
$ws.Range("T5").Value = 0.1120374909379251
